$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "shredding"  (xl/worksheets/sheet1.xml)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A3").Value = 11260
$ws1.Range("B3").Formula = "=0.000046*0.007"
$ws1.Range("C3").Value = 30
$ws1.Range("D3").Formula = "=0.000046*0.007"
$ws1.Range("E3").Value = 0.0136
$ws1.Range("G3").Value = 0.207
$ws1.Range("H3").Formula = "=0.083*0.056"
$ws1.Range("H3").Borders.Item(8).LineStyle = 1
$ws1.Range("I3").Formula = "=0.8*0.056"
$ws1.Range("J3").Formula = "=0.4*0.056"
$ws1.Range("K3").Value = 21

# view state
$ws1.Activate() | Out-Null
$ws1.Range("L5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "extrusion"  (xl/worksheets/sheet2.xml)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A3").Value = 100000
$ws2.Range("B3").Formula = "=4.28*10^-7"
$ws2.Range("C3").Value = 30
$ws2.Range("D3").Formula = "=4.28*10^-7"
$ws2.Range("E3").Value = 0.239
$ws2.Range("G3").Value = 0.207
$ws2.Range("H3").Formula = "=0.083*0.72"
$ws2.Range("H3").Borders.Item(8).LineStyle = 1
$ws2.Range("J3").Formula = "=0.4*0.72"
$ws2.Range("K3").Value = 21

# view state
$ws2.Activate() | Out-Null
$ws2.Range("K9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "granulate"  (xl/worksheets/sheet3.xml) -- row 3 reset to zero
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A3").Value = 0
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 0
$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 0
$ws3.Range("F3").Value = 0
$ws3.Range("G3").Value = 0
$ws3.Range("H3").Value = 0
$ws3.Range("J3").Value = 0
$ws3.Range("K3").Value = 0

# view state
$ws3.Activate() | Out-Null
$ws3.Range("S3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 4: "conditioning"  (xl/worksheets/sheet4.xml) -- becomes active tab
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("Q3").Select()

Write-Output "done"
